$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "42.412.84"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  -1.65%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.532.81"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  -1.52%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "311.12"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -1.77%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "98.33"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +1.79%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.566"
$ws.Range("D7").ClearFormats()
$ws.Range("E8").Value = "  +0.14%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.527"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  -2.76%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "35.58"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -0.40%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0801"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  -1.45%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "7.32"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  -2.13%  "
$ws.Range("E13").Value = "  -0.17%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.920.62"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  -1.56%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "15.77"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  +4.62%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.513.61"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  -2.54%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.833"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  -1.58%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "42.414.86"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  -1.61%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.77"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -1.67%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0₃0946"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -1.63%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.17"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -3.65%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "68.90"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -1.05%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "243.06"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -4.16%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.89"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  -2.07%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.04"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -1.23%  "
$ws.Range("E26").Value = "  +0.00%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "25.97"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -3.24%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.33"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -4.14%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "39.30"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -2.32%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "10.12"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -1.05%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "157.44"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  +2.10%  "
$ws.Range("B32").Value = "Filecoin"
$ws.Range("C32").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "5.69"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -2.34%  "
$ws.Range("B33").Value = "ApeXProtocol"
$ws.Range("C33").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.82"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +15.65%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0794"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  -1.90%  "
$ws.Range("E35").Value = "  -3.27%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.01"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  -5.69%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "18.14"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -4.77%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.16"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -7.61%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.111"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -1.10%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.118"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -0.50%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "4.24"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +9.29%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "21.66"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -3.24%  "
$ws.Range("E43").Value = "  +0.21%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "3.28"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -0.11%  "
$ws.Range("E45").Value = "  -2.90%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.962.29"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -1.75%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "8.93"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -1.39%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.773.69"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -1.60%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "80.73"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -4.21%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.191"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -1.11%  "
$ws.Range("B51").Value = "Aave"
$ws.Range("C51").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "101.29"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -2.97%  "
